$d = $word.ActiveDocument

# Locate the paragraph that ends with:
# "...En caso de habilitar al usuario, todas sus publicaciones pausadas volverán a ser activas."
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*todas sus publicaciones pausadas volverán a ser activas*") {
        $target = $p
        break
    }
}

# Insert the first new bullet paragraph right after the target paragraph.
# InsertParagraphAfter() creates a new paragraph inheriting the same
# paragraph style/numbering (Prrafodelista / numId 1) as $target.
$target.Range.InsertParagraphAfter()
$p1 = $target.Next()
$p1.Range.Text = "Si se deshabilita un rol, ningún usuario podrá acceder con ese rol"

# Insert the second new bullet paragraph after the first one.
$p1.Range.InsertParagraphAfter()
$p2 = $p1.Next()
$p2.Range.Text = "Si se deshabilita una visibilidad no se podrán crear nuevas publicaciones con esa visibilidad, las creadas permanecerán."

# Insert the third new bullet paragraph after the second one.
$p2.Range.InsertParagraphAfter()
$p3 = $p2.Next()
$p3.Range.Text = "Agregamos la funcionalidad de cambiar contraseña para los usuarios, para que puedan cambiar su propia contraseña, el administrador posee la funcionalidad de cambiar su contraseña y la contraseña de cualquier usuario dentro de la funcionalidad de ABM Usuario."

Write-Output "done"
